# v1.1 update reviewer verification status to closed
# close registration wireframe review and verify the updates

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) LH_WF_REGISTRATION_REVIEW: mark the existing review row as closed
# ---------------------------------------------------------------------------
$wsReview = $wb.Worksheets.Item("LH_WF_REGISTRATION_REVIEW")
$wsReview.Range("H2").Value = "closed"
$wsReview.Range("I2").Value = "closed"

# Move the selection/view on this sheet to I2 (without changing which sheet
# is the active tab in the workbook - that stays VERSION-HISTORY).
$wsReview.Activate()
$wsReview.Range("I2").Select()

# ---------------------------------------------------------------------------
# 2) VERSION-HISTORY: add the v1.1 entry that documents closing the review
# ---------------------------------------------------------------------------
$wsHistory = $wb.Worksheets.Item("VERSION-HISTORY")

# Re-use row 2's formatting (fonts, fills, borders, alignment) for row 3 by
# copying formats only, then overwrite the values for the new entry.
$wsHistory.Range("A2:D2").Copy()
$wsHistory.Range("A3:D3").PasteSpecial(-4122)
$wsHistory.Rows.Item(3).RowHeight = 37.5

# The "Updated section" cell for this row uses the slightly different
# (near-white) banding fill instead of a solid white fill.
$wsHistory.Range("C3").Interior.Color = 16382198

$wsHistory.Range("A3").Value = "v1.1"
$wsHistory.Range("B3").Value = "Ahmed Abuzaid"
$wsHistory.Range("C3").Value = "close registration wireframe review and verify the updates"
$wsHistory.Range("D3").Value = 45770

# Make VERSION-HISTORY the active sheet/tab again (matches the original
# workbook state) and move its selection to D10.
$wsHistory.Activate()
$wsHistory.Range("D10").Select()
